# Auto-generated Excel COM-interop script applying the Ixion_Profits.xlsx diff
# to each of the 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 6467.857
$ws.Range("I34").Value = 2235
$ws.Range("J34").Value = 17050
$ws.Range("K34").Value = 2235
$ws.Range("L34").Value = 17050
$ws.Range("M34").Value = -2032
$ws.Range("N34").Value = -17456
$ws.Range("H36").Value = 6467.857
$ws.Range("I36").Value = 2235
$ws.Range("J36").Value = 17050
$ws.Range("K36").Value = 2235
$ws.Range("L36").Value = 17050
$ws.Range("M36").Value = -1520
$ws.Range("N36").Value = -18480
$ws.Range("H137").Value = 1736.5652
$ws.Range("I137").Value = 1472.05
$ws.Range("K137").Value = 4416.15
$ws.Range("M137").Value = -1866.15

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4416.922
$ws.Range("I32").Value = 3226.0195
$ws.Range("J32").Value = 9088.923000000001
$ws.Range("K32").Value = 3226.0195
$ws.Range("L32").Value = 9088.923000000001
$ws.Range("M32").Value = -2939.0195
$ws.Range("N32").Value = -9662.923000000001
$ws.Range("H122").Value = 3210504.5
$ws.Range("I122").Value = 3210504.5
$ws.Range("K122").Value = 9631513.5
$ws.Range("M122").Value = -9629063.5
$ws.Range("H124").Value = 45000
$ws.Range("J124").Value = 45000
$ws.Range("L124").Value = 45000
$ws.Range("N124").Value = -54820
$ws.Range("H132").Value = 2990.0476
$ws.Range("I132").Value = 1417.08
$ws.Range("J132").Value = 5303.2354
$ws.Range("K132").Value = 4251.24
$ws.Range("L132").Value = 15909.7062
$ws.Range("M132").Value = -1721.24
$ws.Range("N132").Value = -20969.7062
$ws.Range("H133").Value = 38950
$ws.Range("J133").Value = 38950
$ws.Range("L133").Value = 38950
$ws.Range("N133").Value = -44010
$ws.Range("H135").Value = 73475
$ws.Range("J135").Value = 73475
$ws.Range("L135").Value = 73475
$ws.Range("N135").Value = -83615

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = $null
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = $null
$ws.Range("N8").Value = $null
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = $null
$ws.Range("H49").Value = 11500
$ws.Range("J49").Value = 11500
$ws.Range("L49").Value = 11500
$ws.Range("N49").Value = -11978
$ws.Range("H86").Value = 55557692
$ws.Range("I86").Value = 111112584
$ws.Range("J86").Value = 2802.3333
$ws.Range("K86").Value = 111112584
$ws.Range("L86").Value = 2802.3333
$ws.Range("M86").Value = -111111461
$ws.Range("N86").Value = -5048.3333
$ws.Range("H89").Value = 55557692
$ws.Range("I89").Value = 111112584
$ws.Range("J89").Value = 2802.3333
$ws.Range("K89").Value = 555562920
$ws.Range("L89").Value = 14011.6665
$ws.Range("M89").Value = -555557304
$ws.Range("N89").Value = -25243.6665
$ws.Range("H134").Value = 2286.8823
$ws.Range("I134").Value = 2047.7273
$ws.Range("J134").Value = 2725.3333
$ws.Range("K134").Value = 6143.1819
$ws.Range("L134").Value = 8175.999899999999
$ws.Range("M134").Value = -3608.1819
$ws.Range("N134").Value = -13245.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 321.17648
$ws.Range("I22").Value = 331.53845
$ws.Range("J22").Value = 287.5
$ws.Range("K22").Value = 331.53845
$ws.Range("L22").Value = 287.5
$ws.Range("M22").Value = 18.46154999999999
$ws.Range("N22").Value = -987.5
$ws.Range("H31").Value = 4945.231
$ws.Range("I31").Value = 1288.24
$ws.Range("J31").Value = 11475.571
$ws.Range("K31").Value = 1288.24
$ws.Range("L31").Value = 11475.571
$ws.Range("M31").Value = -993.24
$ws.Range("N31").Value = -12065.571
$ws.Range("H34").Value = 4945.231
$ws.Range("I34").Value = 1288.24
$ws.Range("J34").Value = 11475.571
$ws.Range("K34").Value = 1288.24
$ws.Range("L34").Value = 11475.571
$ws.Range("M34").Value = -1086.24
$ws.Range("N34").Value = -11879.571
$ws.Range("H58").Value = 1474.4
$ws.Range("I58").Value = 1192.25
$ws.Range("J58").Value = 1712
$ws.Range("K58").Value = 1192.25
$ws.Range("L58").Value = 1712
$ws.Range("M58").Value = -989.25
$ws.Range("N58").Value = -2118
$ws.Range("H99").Value = 15630294
$ws.Range("I99").Value = 2137
$ws.Range("K99").Value = 2137
$ws.Range("M99").Value = -639
$ws.Range("H107").Value = 293.45
$ws.Range("I107").Value = 96.28570999999999
$ws.Range("J107").Value = 399.6154
$ws.Range("K107").Value = 96.28570999999999
$ws.Range("L107").Value = 399.6154
$ws.Range("M107").Value = 1823.71429
$ws.Range("N107").Value = -4239.6154
$ws.Range("H126").Value = 15630294
$ws.Range("I126").Value = 2137
$ws.Range("K126").Value = 6411
$ws.Range("M126").Value = -3941
$ws.Range("H132").Value = 3463
$ws.Range("I132").Value = 3187.3845
$ws.Range("J132").Value = 4179.6
$ws.Range("K132").Value = 9562.1535
$ws.Range("L132").Value = 12538.8
$ws.Range("M132").Value = -7032.1535
$ws.Range("N132").Value = -17598.8
$ws.Range("H136").Value = 1474.4
$ws.Range("I136").Value = 1192.25
$ws.Range("J136").Value = 1712
$ws.Range("K136").Value = 3576.75
$ws.Range("L136").Value = 5136
$ws.Range("M136").Value = -1026.75
$ws.Range("N136").Value = -10236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 23819266
$ws.Range("I117").Value = 50264.5
$ws.Range("J117").Value = 27780766
$ws.Range("K117").Value = 150793.5
$ws.Range("L117").Value = 83342298
$ws.Range("M117").Value = -147351.5
$ws.Range("N117").Value = -83349182
$ws.Range("H121").Value = 904.63794
$ws.Range("J121").Value = 1000.59186
$ws.Range("L121").Value = 3001.77558
$ws.Range("N121").Value = -5621.77558
$ws.Range("H139").Value = 4816.073
$ws.Range("I139").Value = 9693.25
$ws.Range("J139").Value = 2797.9312
$ws.Range("K139").Value = 29079.75
$ws.Range("L139").Value = 8393.793600000001
$ws.Range("M139").Value = -23939.75
$ws.Range("N139").Value = -18673.7936

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4102312
$ws.Range("I122").Value = 2702387
$ws.Range("J122").Value = 12501862
$ws.Range("K122").Value = 8107161
$ws.Range("L122").Value = 37505586
$ws.Range("M122").Value = -8104711
$ws.Range("N122").Value = -37510486
$ws.Range("H132").Value = 4553.963
$ws.Range("I132").Value = 4755.727
$ws.Range("J132").Value = 4415.25
$ws.Range("K132").Value = 14267.181
$ws.Range("L132").Value = 13245.75
$ws.Range("M132").Value = -11737.181
$ws.Range("N132").Value = -18305.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 5499.5
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 5499.5
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 5499.5
$ws.Range("M12").Value = $null
$ws.Range("N12").Value = -5839.5
$ws.Range("H61").Value = 3088.3845
$ws.Range("I61").Value = 2768.625
$ws.Range("K61").Value = 2768.625
$ws.Range("M61").Value = -2566.625
$ws.Range("H113").Value = 3088.3845
$ws.Range("I113").Value = 2768.625
$ws.Range("K113").Value = 2768.625
$ws.Range("M113").Value = -598.625
$ws.Range("H122").Value = 6496350.5
$ws.Range("I122").Value = 7145735.5
$ws.Range("K122").Value = 21437206.5
$ws.Range("M122").Value = -21434756.5
$ws.Range("H132").Value = 13544627
$ws.Range("I132").Value = 14447202
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 43341606
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -43339076
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 15000
$ws.Range("J44").Value = 15000
$ws.Range("L44").Value = 15000
$ws.Range("N44").Value = -16108
$ws.Range("H46").Value = 49425
$ws.Range("J46").Value = 49425
$ws.Range("L46").Value = 49425
$ws.Range("N46").Value = -49887
$ws.Range("H54").Value = 7846.6665
$ws.Range("J54").Value = 7846.6665
$ws.Range("L54").Value = 7846.6665
$ws.Range("N54").Value = -8886.666499999999
$ws.Range("H80").Value = 39550.5
$ws.Range("J80").Value = 39550.5
$ws.Range("L80").Value = 39550.5
$ws.Range("N80").Value = -41546.5
$ws.Range("H83").Value = 39550.5
$ws.Range("J83").Value = 39550.5
$ws.Range("L83").Value = 118651.5
$ws.Range("N83").Value = -128635.5
$ws.Range("H132").Value = 2884.0908
$ws.Range("I132").Value = 2166
$ws.Range("K132").Value = 6498
$ws.Range("M132").Value = -3968
$ws.Range("H134").Value = 49425
$ws.Range("J134").Value = 49425
$ws.Range("L134").Value = 148275
$ws.Range("N134").Value = -153345
$ws.Range("H136").Value = 2328.081
$ws.Range("I136").Value = 2447.4075
$ws.Range("J136").Value = 2005.9
$ws.Range("K136").Value = 6498
$ws.Range("L136").Value = 6017.700000000001
$ws.Range("M136").Value = -4792.2225
$ws.Range("N136").Value = -11117.7

